$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = "Nextion NX4024T032 Display - 3.2"""
$ws.Range("B8").Value = "Nextion NX8048K070 Display - 7"""
$ws.Range("D7").Value = "CASE IS NOT NEEDED, but Amazon is out of stock of just the display as of 2/23/22"

$ws.Range("B7").Hyperlinks.Item(1).TextToDisplay = ""

$ws.Range("D8").Select()
